# Locate the first sentence in the document and append " (Changed main)"
# after it, split across three separate runs (matching the target
# revision's run layout): " (" | "Changed main" | ")".
$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r = $find.Parent
$r.Collapse(0)

# Word normally coalesces newly inserted text into the preceding run when
# the formatting matches exactly. Dropping a transient bookmark at each
# insertion point forces Word to keep the surrounding text split into
# distinct runs, even after the bookmark itself is removed again.

# Run 2: " ("
$seam1 = $d.Bookmarks.Add("seam1", $r)
$r.InsertAfter(" (")
$d.Bookmarks.Item("seam1").Delete() | Out-Null
$r.Collapse(0)

# Run 3: "Changed main"
$seam2 = $d.Bookmarks.Add("seam2", $r)
$r.InsertAfter("Changed main")
$d.Bookmarks.Item("seam2").Delete() | Out-Null
$r.Collapse(0)

# Run 4: ")"
$seam3 = $d.Bookmarks.Add("seam3", $r)
$r.InsertAfter(")")
$d.Bookmarks.Item("seam3").Delete() | Out-Null
$r.Collapse(0)
